{"js": "// Remove all the Lorem-ipsum filler paragraphs except the very first one,\n// keeping the first paragraph and the trailing empty paragraph untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Delete every paragraph except the first and the last (the last is the\n// trailing empty paragraph that closes the document body).\nfor (let i = items.length - 2; i >= 1; i--) {\n  items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove all the Lorem-ipsum filler paragraphs except the very first one,\n# keeping the first paragraph and the trailing empty paragraph untouched.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n\n# Delete paragraphs 2 .. (count-1), walking backwards so indices stay valid\n# as each paragraph is removed. Paragraph 1 (the opening \"Lorem ipsum\u2026\"\n# paragraph) and the final (empty) paragraph are left in place.\nfor ($i = $count - 1; $i -ge 2; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
